$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update changed cell values (rows 2-5) ---
$ws.Range('AA2').Value = -0.7569620253164556
$ws.Range('AB2').Value = 0.06762819835174715
$ws.Range('AC2').Value = -0.8349289569474467
$ws.Range('AD2').Value = 39.28
$ws.Range('AE2').Value = 0
$ws.Range('AF2').Value = 39.28
$ws.Range('AG2').Value = 24
$ws.Range('AH2').Value = 0.1308548204410687
$ws.Range('AI2').Value = 0.7328358208955223
$ws.Range('AJ2').Value = 0.08424008424008425
$ws.Range('AK2').Value = 0.6263048016701461
$ws.Range('AL2').Value = 14.358
$ws.Range('AM2').Value = 13.679
$ws.Range('AN2').Value = -1.557247066286077
$ws.Range('AO2').Value = -2.279774341830338
$ws.Range('AP2').Value = -0.9514747859181731
$ws.Range('AQ2').Value = -2.392938080269025
$ws.Range('B2').Value = '3'
$ws.Range('D2').Value = -0.222
$ws.Range('G2').Value = -0.3682910273410368
$ws.Range('H2').Value = -0.3843057397234325
$ws.Range('I2').Value = -0.258358590642167
$ws.Range('J2').Value = -0.258358590642167
$ws.Range('K2').Value = -69.967
$ws.Range('L2').Value = -0.5522431647407968
$ws.Range('M2').Value = 0
$ws.Range('N2').Value = 0
$ws.Range('O2').Value = -0
$ws.Range('S2').Value = 0
$ws.Range('U2').Value = 15.28
$ws.Range('V2').Value = 0.05856650057493294
$ws.Range('W2').Value = -2.569811320754717
$ws.Range('X2').Value = 0.07204832365549309
$ws.Range('Y2').Value = -2.64185964441021
$ws.Range('Z2').Value = 2.540525365951474
$ws.Range('AA3').Value = -0.1918997107039536
$ws.Range('AB3').Value = 0.06695005156716119
$ws.Range('AC3').Value = -0.2588497622711148
$ws.Range('AD3').Value = 10.8
$ws.Range('AE3').Value = 0
$ws.Range('AF3').Value = 10.8
$ws.Range('AG3').Value = 4.94
$ws.Range('AH3').Value = 0.08933002481389579
$ws.Range('AI3').Value = 0.1803005008347245
$ws.Range('AJ3').Value = 0.04294158553546593
$ws.Range('AK3').Value = 0.0914137675795707
$ws.Range('AL3').Value = 1.31
$ws.Range('AM3').Value = 0.783
$ws.Range('AN3').Value = 2.245322245322245
$ws.Range('AO3').Value = -1.519083969465649
$ws.Range('AP3').Value = 1.027027027027027
$ws.Range('AQ3').Value = -2.541507024265645
$ws.Range('G3').Value = 0.01520417028670721
$ws.Range('H3').Value = 0.01520417028670721
$ws.Range('I3').Value = -0.01728931364031277
$ws.Range('J3').Value = -0.01728931364031277
$ws.Range('K3').Value = -0.147
$ws.Range('L3').Value = -0.001277150304083406
$ws.Range('U3').Value = 5.86
$ws.Range('V3').Value = 0.05322434150772026
$ws.Range('W3').Value = -0.00310126582278481
$ws.Range('X3').Value = 0.07204832365549309
$ws.Range('Y3').Value = -0.07514958947827789
$ws.Range('Z3').Value = 11.09932497589199
$ws.Range('AA4').Value = -0.7569620253164556
$ws.Range('AB4').Value = 0.07796693163099112
$ws.Range('AC4').Value = -0.8349289569474467
$ws.Range('AD4').Value = 25.7
$ws.Range('AF4').Value = 25.7
$ws.Range('AG4').Value = 18.65
$ws.Range('AH4').Value = 0.2827282728272827
$ws.Range('AI4').Value = -1.992248062015503
$ws.Range('AJ4').Value = 0.2224209898628503
$ws.Range('AK4').Value = -0.9348370927318294
$ws.Range('AL4').Value = 12.8
$ws.Range('AM4').Value = 12.648
$ws.Range('AN4').Value = -0.8741496598639455
$ws.Range('AO4').Value = -2.3359375
$ws.Range('AP4').Value = -0.6343537414965986
$ws.Range('AQ4').Value = -2.364010120177103
$ws.Range('D4').Value = -0.222
$ws.Range('G4').Value = -95.28225806451613
$ws.Range('H4').Value = -98.99193548387098
$ws.Range('I4').Value = -60.28225806451613
$ws.Range('J4').Value = -60.28225806451613
$ws.Range('K4').Value = -68.09999999999999
$ws.Range('L4').Value = -137.2983870967742
$ws.Range('U4').Value = 7.05
$ws.Range('V4').Value = 0.1081288343558282
$ws.Range('W4').Value = -2.569811320754717
$ws.Range('X4').Value = 0.08550957876161185
$ws.Range('Y4').Value = -2.655320899516329
$ws.Range('Z4').Value = 0.01255696202531646
$ws.Range('AB5').Value = 0.06762819835174715
$ws.Range('AD5').Value = 2.78
$ws.Range('AE5').Value = 0
$ws.Range('AF5').Value = 2.78
$ws.Range('AG5').Value = 0.4099999999999997
$ws.Range('AH5').Value = 0.0314550803349174
$ws.Range('AI5').Value = 0.4212121212121212
$ws.Range('AJ5').Value = 0.004766887571212647
$ws.Range('AK5').Value = 0.09692671394799049
$ws.Range('AL5').Value = 0.248
$ws.Range('AM5').Value = 0.248
$ws.Range('AN5').Value = -4.384858044164037
$ws.Range('AO5').Value = -3.399193548387097
$ws.Range('AP5').Value = -0.6466876971608828
$ws.Range('AQ5').Value = -3.399193548387097
$ws.Range('B5').Value = 'Xplora Technologies AS (OB:XPLRA)'
$ws.Range('G5').Value = -0.1036936936936937
$ws.Range('H5').Value = -0.1207207207207207
$ws.Range('I5').Value = -0.07594594594594595
$ws.Range('J5').Value = -0.07594594594594595
$ws.Range('K5').Value = -1.72
$ws.Range('L5').Value = -0.154954954954955
$ws.Range('U5').Value = 2.37
$ws.Range('V5').Value = 0.02768691588785047
$ws.Range('X5').Value = 0.06906508806457258

# --- Clear individually-removed cells ---
$ws.Range('T2').ClearContents()
$ws.Range('AA5').ClearContents()
$ws.Range('AC5').ClearContents()
$ws.Range('W5').ClearContents()
$ws.Range('Y5').ClearContents()
$ws.Range('Z5').ClearContents()

# --- Remove row 6 entirely (Huddly AS row dropped; rows shift up / dimension shrinks) ---
$ws.Rows(6).Delete()
